# Add ability to supply row to read data from:
# populate a new data row (row 3) under the existing header row, matching
# the same column layout used by the header row (A..AD), then move the
# selection onto the newly entered data (AC3) as the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowIndex = 3
$values = @(3, 4, 2, 6, 7, 3, 2, 4, 5, 6, 23, 2, 45, 7, 8, 4, 32, 5, 6, 4, 32, 5, 2, 54, 2, 3, 2, 7, 3, 2)

for ($col = 1; $col -le $values.Length; $col++) {
    $ws.Cells.Item($rowIndex, $col).Value = $values[$col - 1]
}

$ws.Range("AC3").Select()
